$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet updates
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# The note explaining the "electrolysis with guaranteed clean electricity"
# row is replaced with a note about the BAU production shares / IRA repeal.
$about.Range("A12").Value = "the BAU production shares, representing a repeal of IRA tax credits"

# ---------------------------------------------------------------------------
# "RHPF" sheet updates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("RHPF")

# Rename the "hydrocarbon partial oxidation" pathway to
# "thermochemical water splitting" (header row + row label).
$ws.Range("F1").Value = "thermochemical water splitting"
$ws.Range("A6").Value = "thermochemical water splitting"

# Update the fraction values for the "electrolysis" row (row 2):
# from 0 across the board to 0.05.
$ws.Range("B2:H2").Value = 0.05

# Update the fraction values for the "natural gas reforming" row (row 3):
# from 0 across the board to 0.95.
$ws.Range("B3:H3").Value = 0.95

# Update the fraction values for the "electrolysis with guaranteed clean
# electricity" row (row 7): from 1 across the board to 0.
$ws.Range("B7:H7").Value = 0

# Columns G:H did not previously carry the right-aligned numeric style that
# columns B:F get from the column default; match that formatting now that
# those cells hold real data, mirroring columns B:F (style index 3).
$ws.Range("G2:H3").HorizontalAlignment = -4152
$ws.Range("G7:H7").HorizontalAlignment = -4152

# Move the active-cell selection on the RHPF sheet (cosmetic, matches diff).
$ws.Range("B2:H3").Select() | Out-Null

# The "About" sheet remains the active/visible tab, so re-select its
# active cell last (cosmetic, matches diff).
$about.Range("B13").Select() | Out-Null
